# Clean up versions of code and manuscript
# - Rename the trailing "Sheet1" tab to "exo_wage"
# - Re-point that sheet's view (zoom 110%, selection on L17)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "exo_wage"

# Make it the active sheet/window so the Zoom + selection land on it
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$ws.Range("L17").Select()
